# Applies the cryptos-list refresh described by the commit diff.
# Numeric-looking "Price" strings (e.g. "604.83") are prefixed with a
# leading apostrophe so Excel stores them as text (matching the original
# t="inlineStr" cells) instead of auto-converting them to numbers; the
# Style reset afterwards clears the quote-prefix formatting flag Excel
# adds for that trick, so no unrelated styling changes remain.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.522.52'
$ws.Range('E2').Value = '  +1.19%  '
$ws.Range('D3').Value = '3.164.22'
$ws.Range('E3').Value = '  +0.73%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range("D5").Value = "`'604.83"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.71%  '
$ws.Range("D6").Value = "`'144.82"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.87%  '
$ws.Range("D7").Value = "`'1.01"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.70%  '
$ws.Range('D8').Value = '3.152.88'
$ws.Range('E8').Value = '  +0.65%  '
$ws.Range("D9").Value = "`'0.524"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.46%  '
$ws.Range("D10").Value = "`'0.150"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.25%  '
$ws.Range("D11").Value = "`'5.40"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.80%  '
$ws.Range("D12").Value = "`'0.473"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.22%  '
$ws.Range("D13").Value = "`'0.0000255"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.81%  '
$ws.Range("D14").Value = "`'35.64"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.45%  '
$ws.Range('D15').Value = '3.682.31'
$ws.Range('E15').Value = '  +0.93%  '
$ws.Range('E16').Value = '  +2.37%  '
$ws.Range('D17').Value = '64.664.06'
$ws.Range('E17').Value = '  +1.30%  '
$ws.Range('D18').Value = '3.169.11'
$ws.Range('E18').Value = '  +1.00%  '
$ws.Range("D19").Value = "`'6.90"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.87%  '
$ws.Range("D20").Value = "`'483.33"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.21%  '
$ws.Range("D21").Value = "`'14.71"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.57%  '
$ws.Range("D22").Value = "`'0.715"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range("D23").Value = "`'7.72"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.42%  '
$ws.Range("D24").Value = "`'87.81"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.88%  '
$ws.Range("D25").Value = "`'13.53"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.44%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range("D27").Value = "`'2.78"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.81%  '
$ws.Range("D28").Value = "`'8.45"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.58%  '
$ws.Range("D29").Value = "`'7.19"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.70%  '
$ws.Range('E30').Value = '  +0.46%  '
$ws.Range('E31').Value = '  -8.69%  '
$ws.Range("D32").Value = "`'27.14"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.47%  '
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range("D34").Value = "`'2.70"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.11%  '
$ws.Range('E35').Value = '  -1.32%  '
$ws.Range("D36").Value = "`'6.07"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.04%  '
$ws.Range('D37').Value = '0.0₃0764'
$ws.Range('E37').Value = '  -3.81%  '
$ws.Range("D38").Value = "`'53.01"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.25%  '
$ws.Range("D39").Value = "`'3.04"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.18%  '
$ws.Range("D40").Value = "`'443.54"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.37%  '
$ws.Range('E41').Value = '  +0.71%  '
$ws.Range("D42").Value = "`'0.119"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.41%  '
$ws.Range("D43").Value = "`'8.30"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.58%  '
$ws.Range('D44').Value = '2.873.58'
$ws.Range('E44').Value = '  +0.89%  '
$ws.Range("D45").Value = "`'0.264"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.51%  '
$ws.Range("D46").Value = "`'2.49"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.12%  '
$ws.Range("D47").Value = "`'2.25"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.50%  '
$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D48").Value = "`'0.998"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = "`'26.18"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.62%  '
$ws.Range('E50').Value = '  +0.51%  '
$ws.Range("D51").Value = "`'121.82"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.38%  '
